$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Obrigatorio" column (E) from "N" to "S" for rows 2 through 8
$ws.Range("E2:E8").Value = "S"
